$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs"/"balls"/"fours" columns store numeric-looking values as TEXT.
# Force a Text number format before assigning so Excel keeps them as text
# (t="str"/shared-string text) instead of auto-converting to numbers.

$cells = @{
    "C2" = "2"
    "D2" = "4"
    "C3" = "1"
    "D3" = "2"
    "C4" = "0"
    "C5" = "9"
    "D5" = "6"
    "E5" = "1"
    "C6" = "0"
    "D6" = "1"
    "E6" = "0"
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}
